# Doing Updates for Financials
# Update the Total Revenue (row 9) and Cost of Revenue (row 10) figures
# on the BDCO income-statement sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total Revenue (row 9)
$ws.Range("D9").Value = 249300
$ws.Range("E9").Value = 173800
$ws.Range("F9").Value = 204900
$ws.Range("G9").Value = 370900
$ws.Range("H9").Value = 409800
$ws.Range("I9").Value = 350700
$ws.Range("J9").Value = 1200

# Cost of Revenue (row 10)
$ws.Range("D10").Value = 9200
$ws.Range("E10").Value = -5900
$ws.Range("F10").Value = 16800
$ws.Range("G10").Value = 17800
$ws.Range("H10").Value = -300
$ws.Range("I10").Value = 1400
$ws.Range("J10").Value = 1100
